$wb = $excel.ActiveWorkbook

# Step 1: clear all existing content on every worksheet so the shared-string
# table is rebuilt from scratch in the same order a fresh export would use.
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.ClearContents()
}

# --- Funciones_Objetivo ---
$ws = $wb.Worksheets.Item("Funciones_Objetivo")
$ws.Range("A1").Formula = "=""Leader_Expr"""
$ws.Range("B1").Formula = "=""Follower_Expr"""
$rng = $ws.Range("A1:B1")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A2").Formula = "=""-x - 3y_1 + 2y_2"""
$ws.Range("B2").Formula = "=""-y_1"""
$rng = $ws.Range("A2:B2")
$rng.Copy()
$rng.PasteSpecial(-4163)

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
$ws.Range("A1").Formula = "=""Expression"""
$ws.Range("B1").Formula = "=""Function_Evaluation"""
$ws.Range("C1").Formula = "=""Restriction_Set_Type"""
$ws.Range("D1").Formula = "=""MIU_value"""
$rng = $ws.Range("A1:D1")
$rng.Copy()
$rng.PasteSpecial(-4163)

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A1").Formula = "=""Expression"""
$ws.Range("B1").Formula = "=""Function_Evaluation"""
$ws.Range("C1").Formula = "=""Restriction_Set_Type"""
$ws.Range("D1").Formula = "=""Lambda_value"""
$ws.Range("E1").Formula = "=""Beta_value"""
$ws.Range("F1").Formula = "=""Gamma_value"""
$rng = $ws.Range("A1:F1")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A2").Formula = "=""7.35 - y_1"""
$ws.Range("B2").Formula = "=""-7.35"""
$ws.Range("C2").Formula = "=""J_0_L0_v"""
$ws.Range("D2").Formula = "=""0.03"""
$ws.Range("E2").Formula = "=""0.4"""
$ws.Range("F2").Formula = "=""0"""
$rng = $ws.Range("A2:F2")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A3").Formula = "=""-7.35 + y_1"""
$ws.Range("B3").Formula = "=""3.3499999999999996"""
$ws.Range("C3").Formula = "=""J_0_L0_v"""
$ws.Range("D3").Formula = "=""0.85"""
$ws.Range("E3").Formula = "=""5.0"""
$ws.Range("F3").Formula = "=""5.2"""
$rng = $ws.Range("A3:F3")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A4").Formula = "=""-3.7499999999999964 - 2x + y_1 + 4y_2"""
$ws.Range("B4").Formula = "=""-12.250000000000004"""
$ws.Range("C4").Formula = "=""J_0_LP_v"""
$ws.Range("D4").Formula = "=""0.73"""
$ws.Range("E4").Formula = "=""0"""
$ws.Range("F4").Formula = "=""7.4"""
$rng = $ws.Range("A4:F4")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A5").Formula = "=""-67.42 + 8x + y_1"""
$ws.Range("B5").Formula = "=""18.550000000000004"""
$ws.Range("C5").Formula = "=""J_Ne_L0_v"""
$ws.Range("D5").Formula = "=""0.02"""
$ws.Range("E5").Formula = "=""5.8"""
$ws.Range("F5").Formula = "=""9.5"""
$rng = $ws.Range("A5:F5")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A6").Formula = "=""-5.5 - 2x - 2y_1"""
$ws.Range("B6").Formula = "=""-17.5"""
$ws.Range("C6").Formula = "=""J_Ne_L0_v"""
$ws.Range("D6").Formula = "=""0.72"""
$ws.Range("E6").Formula = "=""2.9"""
$ws.Range("F6").Formula = "=""3.2"""
$rng = $ws.Range("A6:F6")
$rng.Copy()
$rng.PasteSpecial(-4163)

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A1").Formula = "=""x"""
$ws.Range("B1").Formula = "=""y_1"""
$ws.Range("C1").Formula = "=""y_2"""
$rng = $ws.Range("A1:C1")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A2").Formula = "=""7.4"""
$ws.Range("B2").Formula = "=""7.35"""
$ws.Range("C2").Formula = "=""2.8"""
$rng = $ws.Range("A2:C2")
$rng.Copy()
$rng.PasteSpecial(-4163)

# --- Vector_bf ---
$ws = $wb.Worksheets.Item("Vector_bf")
$ws.Range("A1").Formula = "=""vec_bf"""
$rng = $ws.Range("A1")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A2").Formula = "=""0.8700000000000001"""
$rng = $ws.Range("A2")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A3").Formula = "=""-2.92"""
$rng = $ws.Range("A3")
$rng.Copy()
$rng.PasteSpecial(-4163)

# --- Vector_BF ---
$ws = $wb.Worksheets.Item("Vector_BF")
$ws.Range("A1").Formula = "=""vec_BF"""
$rng = $ws.Range("A1")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A2").Formula = "=""-39.6"""
$rng = $ws.Range("A2")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A3").Formula = "=""-1.5999999999999988"""
$rng = $ws.Range("A3")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A4").Formula = "=""-2.0"""
$rng = $ws.Range("A4")
$rng.Copy()
$rng.PasteSpecial(-4163)

# --- Vector_Alpha ---
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A1").Formula = "=""vec_alpha"""
$rng = $ws.Range("A1")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A2").Value = 0.0
$ws.Range("A3").Value = 0.0

